$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 33: Wages Expense - top-up the Debit formula ---
$ws.Range("D33").Formula = "=60000+260000"

# --- Row 34: TRANSFER BCA - add the extra 5,070,000 transfer ---
$ws.Range("D34").Formula = "=2216000+154800000+1130000+55000+5070000"

# --- Row 35: A/R - add the extra 14,262,000 receivable ---
$ws.Range("C35").Formula = "=154800000+14262000"

# --- Row 36: new entry - BELI kresek ---
$ws.Range("B36").Value = "BELI kresek"
$ws.Range("D36").Value = 82000

# --- Row 37: new entry - SALES - cash/retail ---
$ws.Range("B37").Value = "SALES - cash/retail"
$ws.Range("C37").Formula = "=8493475+11686525-14262000"

# --- Row 38: new entry - SELISIH - lebih ---
$ws.Range("B38").Value = "SELISIH - lebih"
$ws.Range("C38").Value = 60000

# --- Row 39: new entry - SETOR KE BANK ---
$ws.Range("B39").Value = "SETOR KE BANK"
$ws.Range("D39").Value = 11000000

# --- Row 40: new day, 6-Feb-2021 (serial 44233) - Wages Expense ---
$ws.Range("A40").Value = 44233
$ws.Range("B40").Value = "Wages Expense"
$ws.Range("D40").Formula = "=60000"

# --- Row 41: new entry - A/R ---
$ws.Range("B41").Value = "A/R"
$ws.Range("C41").Formula = "=500000+16610000+12250000"

# --- Row 42: new entry - TRANSFER BCA ---
$ws.Range("B42").Value = "TRANSFER BCA"
$ws.Range("D42").Formula = "=28860000+2100000"

# --- Row 43: new entry - BENSIN - rush ---
$ws.Range("B43").Value = "BENSIN - rush"
$ws.Range("D43").Value = 250000

# --- View state: scroll the frozen pane down and move the active selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 39
$ws.Range("C60").Select()
